$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final data table for the Fgf1-Fgfr3 LR-pair sheet, following Dr Hou's
# revised advice: the sending-cluster set now includes "M2" in addition to
# ECs / FAPs / sCs, and counts/weights are recomputed accordingly.
$data = @(
  @("ECs", "Fgf1", "Fgfr3", "ECs", 3, 1, 0.5247423333333333, 1.574227, 0.044891155074209, 0.044891155074209, 3, 1, 6.408252, 19.224756, 0.8583439096634812, 0.8583439096634812, 3.362681107068, 30.264129963612, 0.03853204955570617, 0.03853204955570617),
  @("ECs", "Fgf1", "Fgfr3", "FAPs", 3, 1, 0.5247423333333333, 1.574227, 0.044891155074209, 0.044891155074209, 3, 1, 0.6824433333333334, 2.04733, 0.09140887075868921, 0.09140887075868921, 0.3581069071011111, 3.22296216391, 0.004103449792386646, 0.004103449792386646),
  @("ECs", "Fgf1", "Fgfr3", "sCs", 3, 1, 0.5247423333333333, 1.574227, 0.044891155074209, 0.044891155074209, 3, 1, 0.3751373333333333, 1.125412, 0.05024721957782962, 0.05024721957782963, 0.1968504396137777, 1.771653956524, 0.00225565572611618, 0.00225565572611618),
  @("FAPs", "Fgf1", "Fgfr3", "ECs", 3, 1, 4.531154666666667, 13.593464, 0.3876355191593572, 0.3876355191593572, 3, 1, 6.408252, 19.224756, 0.8583439096634812, 0.8583439096634812, 29.036780954976, 261.331028594784, 0.3327245870396759, 0.3327245870396759),
  @("FAPs", "Fgf1", "Fgfr3", "FAPs", 3, 1, 4.531154666666667, 13.593464, 0.3876355191593572, 0.3876355191593572, 3, 1, 0.6824433333333334, 2.04733, 0.09140887075868921, 0.09140887075868921, 3.092256294568889, 27.83030665112, 0.03543332507231508, 0.03543332507231508),
  @("FAPs", "Fgf1", "Fgfr3", "sCs", 3, 1, 4.531154666666667, 13.593464, 0.3876355191593572, 0.3876355191593572, 3, 1, 0.3751373333333333, 1.125412, 0.05024721957782962, 0.05024721957782963, 1.699805278574222, 15.298247507168, 0.0194776070473662, 0.0194776070473662),
  @("M2", "Fgf1", "Fgfr3", "ECs", 1, 0.3333333333333333, 0.08000233333333333, 0.240007, 0.006844115528380393, 0.006844115528380393, 3, 1, 6.408252, 19.224756, 0.8583439096634812, 0.8583439096634812, 0.512675112588, 4.614076013292, 0.005874604880818569, 0.005874604880818569),
  @("M2", "Fgf1", "Fgfr3", "FAPs", 1, 0.3333333333333333, 0.08000233333333333, 0.240007, 0.006844115528380393, 0.006844115528380393, 3, 1, 0.6824433333333334, 2.04733, 0.09140887075868921, 0.09140887075868921, 0.05459705903444444, 0.49137353131, 0.0006256128717912613, 0.0006256128717912613),
  @("M2", "Fgf1", "Fgfr3", "sCs", 1, 0.3333333333333333, 0.08000233333333333, 0.240007, 0.006844115528380393, 0.006844115528380393, 3, 1, 0.3751373333333333, 1.125412, 0.05024721957782962, 0.05024721957782963, 0.03001186198711111, 0.270106757884, 0.000343897775770563, 0.000343897775770563),
  @("sCs", "Fgf1", "Fgfr3", "ECs", 3, 1, 6.553315, 19.659945, 0.5606292102380533, 0.5606292102380533, 3, 1, 6.408252, 19.224756, 0.8583439096634812, 0.8583439096634812, 41.99529395538001, 377.95764559842, 0.4812126681872805, 0.4812126681872805),
  @("sCs", "Fgf1", "Fgfr3", "FAPs", 3, 1, 6.553315, 19.659945, 0.5606292102380533, 0.5606292102380533, 3, 1, 0.6824433333333334, 2.04733, 0.09140887075868921, 0.09140887075868921, 4.472266132983334, 40.25039519685, 0.05124648302219622, 0.05124648302219622),
  @("sCs", "Fgf1", "Fgfr3", "sCs", 3, 1, 6.553315, 19.659945, 0.5606292102380533, 0.5606292102380533, 3, 1, 0.3751373333333333, 1.125412, 0.05024721957782962, 0.05024721957782963, 2.458393113593333, 22.12553802234, 0.02817005902857667, 0.02817005902857668)
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $row = $data[$i]
  $r = $i + 2
  for ($c = 0; $c -lt $row.Length; $c++) {
    $ws.Cells.Item($r, $c + 1).Value = $row[$c]
  }
}
